$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.062.72"
$ws.Range("E2").Value = "  -6.19%  "
$ws.Range("D3").Value = "2.433.34"
$ws.Range("E3").Value = "  -9.28%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'528.35"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").Value = "'146.64"
$ws.Range("E6").Value = "  -6.94%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("D9").Value = "'0.0982"
$ws.Range("E9").Value = "  -6.85%  "
$ws.Range("D11").Value = "'5.32"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("E12").Value = "  -5.77%  "
$ws.Range("D13").Value = "2.867.30"
$ws.Range("E13").Value = "  -9.21%  "
$ws.Range("D14").Value = "'23.94"
$ws.Range("E14").Value = "  -7.96%  "
$ws.Range("D15").Value = "59.002.50"
$ws.Range("E15").Value = "  -6.03%  "
$ws.Range("E16").Value = "  -6.79%  "
$ws.Range("D17").Value = "2.484.55"
$ws.Range("E17").Value = "  -7.45%  "
$ws.Range("E18").Value = "  -7.82%  "
$ws.Range("D19").Value = "'4.28"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("D20").Value = "'322.11"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("D21").Value = "'0.967"
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("E22").Value = "  -9.62%  "
$ws.Range("E23").Value = "  -7.74%  "
$ws.Range("D24").Value = "'60.17"
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("D25").Value = "'0.160"
$ws.Range("E25").Value = "  -4.22%  "
$ws.Range("D26").Value = "'0.970"
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("E27").Value = "  -5.84%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.28"
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'6.74"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").Value = "  -6.24%  "
$ws.Range("D31").Value = "0.0₃0768"
$ws.Range("E31").Value = "  -9.88%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'157.97"
$ws.Range("E33").Value = "  -4.46%  "
$ws.Range("D34").Value = "'4.48"
$ws.Range("E34").Value = "  -6.41%  "
$ws.Range("D35").Value = "'18.24"
$ws.Range("E35").Value = "  -6.57%  "
$ws.Range("E36").Value = "  -6.77%  "
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'5.73"
$ws.Range("E38").Value = "  -7.14%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'310.08"
$ws.Range("E39").Value = "  -8.53%  "
$ws.Range("D40").Value = "'36.60"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("D41").Value = "'0.841"
$ws.Range("E41").Value = "  -9.86%  "
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = "  -6.21%  "
$ws.Range("D43").Value = "'0.995"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'10.67"
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.576"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0928"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("D47").Value = "'0.0519"
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("D48").Value = "'18.85"
$ws.Range("E48").Value = "  -9.13%  "
$ws.Range("D49").Value = "1.972.94"
$ws.Range("E49").Value = "  -5.48%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'18.35"
$ws.Range("E50").Value = "  -9.76%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0226"
$ws.Range("E51").Value = "  -5.45%  "
